# Fix hard-coded UnitMass (column C) values that were overriding the
# strings grabbed from the UI text fields in main_test.
$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Sheet1")

$updates = @{
    "C3"  = 55
    "C5"  = 29
    "C6"  = 71
    "C8"  = 115
    "C9"  = 40
    "C10" = 175
    "C11" = 67
    "C12" = 27
    "C13" = 56
    "C14" = 103
    "C15" = 231
    "C16" = 51
    "C17" = 79
    "C18" = 77
    "C19" = 57
    "C20" = 32
    "C21" = 128
    "C23" = 39
    "C25" = 30
    "C26" = 58
    "C28" = 31
    "C29" = 82
    "C30" = 83
    "C31" = 60
    "C32" = 98
    "C33" = 44
    "C34" = 42
    "C35" = 19
    "C36" = 84
    "C37" = 68
    "C38" = 18
    "C39" = 85
    "C40" = 70
    "C41" = 74
    "C42" = 91
}

foreach ($addr in $updates.Keys) {
    $ws.Range($addr).Value = $updates[$addr]
}
